# Update of Excel Modules Files
# Fill in the (previously empty) "Name" / "ChefModule" / "ElementName1" / "ElementName2"
# columns for each module row, drop the now-unused "ElementName3" column (F), and
# refresh the best-fit column widths / selection the same way a user would after
# typing the new data in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Code (already present, kept for reference) | Name | ChefModule | ElementName1 | ElementName2
$rows = @(
    ,@("G3EI131", "pede. Suspendisse dui.",       "EL Haddad",        "Nullam feugiat placerat",   "varius et, euismod")
    ,@("G3EI132", "a nunc. In",                    "Badir",            "sodales nisi magna",        "elementum sem, vitae")
    ,@("G3EI133", "amet metus. Aliquam",           "Ezzine",           "Cras vulputate velit",      "scelerisque neque sed")
    ,@("G3EI134", "quam vel sapien",                "El Alami Hassoun", "Nunc mauris elit,",         "libero et tristique")
    ,@("G3EI135", "feugiat nec, diam.",             "Lazaar",           "pellentesque. Sed dictum.", "ridiculus mus. Proin")
    ,@("G3EI136", "nonummy. Fusce fermentum",       "El Haddad",        "neque pellentesque massa",  "Mauris eu turpis.")
    ,@("G3EI141", "a, arcu. Sed",                   "EL Haddad",        "sit amet risus.",           "Nulla facilisi. Sed")
    ,@("G3EI142", "Suspendisse eleifend. Cras",     "El Alami Hassoun", "velit dui, semper",         "ligula elit, pretium")
    ,@("G3EI143", "ante. Nunc mauris",              "Badir",            "tortor at risus.",          "felis. Donec tempor,")
    ,@("G3EI144", "lobortis quam a",                "Ezzine",           "euismod est arcu",          "ligula eu enim.")
    ,@("G3EI145", "rhoncus. Nullam velit",          "Ben Achrab",       "ut dolor dapibus",          "commodo tincidunt nibh.")
    ,@("G3EI146", "Donec tincidunt. Donec",         "EL Haddad",        "ornare tortor at",          "ac, feugiat non,")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# The "ElementName3" column (F) is no longer used - clear it out entirely.
$ws.Range("F1:F13").ClearContents()

# Let the newly-filled columns re-compute their best-fit width, like Excel does
# automatically when a bestFit column's contents change. Column F keeps its
# original width since it no longer holds any data.
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()

# Leave the selection where the user finished editing.
[void]$ws.Range("F9").Select()
